# Apply updated cryptocurrency data to the worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '45.141.25'
$ws.Range("E2").Value = '  +1.06%  '
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.359.55'
$ws.Range("E3").Value = '  -0.78%  '
# Row 4
$ws.Range("E4").Value = '  +0.23%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.95'
$ws.Range("E5").Value = '  -0.99%  '
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '107.59'
$ws.Range("E6").Value = '  -1.40%  '
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.629'
$ws.Range("E7").Value = '  -0.37%  '
# Row 8
$ws.Range("E8").Value = '  +0.18%  '
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.609'
$ws.Range("E9").Value = '  -3.18%  '
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.76'
$ws.Range("E10").Value = '  -2.89%  '
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0917'
$ws.Range("E11").Value = '  -0.90%  '
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.49'
$ws.Range("E12").Value = '  -1.96%  '
# Row 13
$ws.Range("E13").Value = '  +0.65%  '
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.975'
$ws.Range("E14").Value = '  -4.12%  '
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.718.53'
$ws.Range("E15").Value = '  -0.20%  '
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.25'
$ws.Range("E16").Value = '  -2.55%  '
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.354.88'
$ws.Range("E17").Value = '  -1.34%  '
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '45.071.04'
$ws.Range("E18").Value = '  +0.80%  '
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.12'
$ws.Range("E19").Value = '  +10.93%  '
# Row 20
$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.24'
$ws.Range("E20").Value = '  -5.31%  '
# Row 21
$ws.Range("B21").Value = 'ShibaInu'
$ws.Range("C21").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0000106'
$ws.Range("E21").Value = '  -1.09%  '
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.16'
$ws.Range("E22").Value = '  -2.56%  '
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.55'
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '258.22'
$ws.Range("E24").Value = '  -4.18%  '
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.36'
$ws.Range("E25").Value = '  +3.19%  '
# Row 26
$ws.Range("E26").Value = '  +0.17%  '
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.03'
$ws.Range("E27").Value = '  -1.85%  '
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.22'
$ws.Range("E28").Value = '  -4.06%  '
# Row 29
$ws.Range("E29").Value = '  +0.00%  '
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0982'
$ws.Range("E30").Value = '  +7.06%  '
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '22.19'
$ws.Range("E31").Value = '  -3.01%  '
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '37.26'
$ws.Range("E32").Value = '  -5.10%  '
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '166.92'
$ws.Range("E33").Value = '  -1.53%  '
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.00'
$ws.Range("E34").Value = '  +4.45%  '
# Row 35
$ws.Range("E35").Value = '  -1.58%  '
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.117'
$ws.Range("E36").Value = '  -1.60%  '
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.71'
$ws.Range("E37").Value = '  -1.29%  '
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.99'
$ws.Range("E38").Value = '  +5.14%  '
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0352'
$ws.Range("E39").Value = '  -3.33%  '
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.89'
$ws.Range("E40").Value = '  -1.81%  '
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.78'
$ws.Range("E41").Value = '  +1.23%  '
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '97.93'
$ws.Range("E42").Value = '  -6.77%  '
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '69.76'
$ws.Range("E43").Value = '  -2.61%  '
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.229'
$ws.Range("E44").Value = '  -4.58%  '
# Row 45
$ws.Range("B45").Value = 'Celestia'
$ws.Range("C45").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '12.89'
$ws.Range("E45").Value = '  -6.14%  '
# Row 46
$ws.Range("B46").Value = 'FirstDigitalUSD'
$ws.Range("C46").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.00'
$ws.Range("E46").Value = '  +0.06%  '
# Row 47
$ws.Range("B47").Value = 'Maker'
$ws.Range("C47").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.833.14'
$ws.Range("E47").Value = '  +9.75%  '
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '83.90'
$ws.Range("E48").Value = '  +5.76%  '
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.73'
$ws.Range("E49").Value = '  +5.32%  '
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '111.24'
$ws.Range("E50").Value = '  -4.79%  '
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '9.16'
$ws.Range("E51").Value = '  +1.35%  '
